$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")

# Update the "custo" value for the electronics row; this value feeds the
# "Planilha1!D4-C2" formula on Planilha2 and the "Planilha2!D2" formula on
# Planilha1, so both recalculate automatically.
$ws1.Range("D4").Value = 520.0

# Add the new "VALOR" / "QUANTIDADE" header row used by the new example
# calculation below it.
$ws1.Range("D5").Value = "VALOR"
$ws1.Range("E5").Value = "QUANTIDADE"

# Add sample value/quantity data together with a sum formula and a
# multiplication formula so that, after the user edits D6/E6, F6 and G6
# keep themselves up to date automatically.
$ws1.Range("D6").Value = 500.3
$ws1.Range("E6").Value = 2.0
$ws1.Range("F6").Formula = "=SUM(D6:E6)"
$ws1.Range("G6").Formula = "=D6*E6"

# Make sure row 8 exists (matches the sheet growing by one row) and give
# rows 6-8 the tighter row height used for the new data rows.
$ws1.Rows.Item(6).RowHeight = 13.8
$ws1.Rows.Item(7).RowHeight = 13.8
$ws1.Rows.Item(8).RowHeight = 13.8

# Leave the selection on the newly added header cell.
$ws1.Range("D5").Select()
